$wb = $excel.ActiveWorkbook
$wsPOP = $wb.Worksheets.Item("POP")
$wsPest = $wb.Worksheets.Item("Pesticide")

# --- POP sheet updates ---
    $wsPOP.Range("B2").Value = 1.5883
    $wsPOP.Range("C2").Value = 0.5841
    $wsPOP.Range("D2").Value = 0.0065
    $wsPOP.Range("B3").Value = -0.1657
    $wsPOP.Range("C3").Value = 0.3853
    $wsPOP.Range("D3").Value = 0.6672
    $wsPOP.Range("B4").Value = 0.288
    $wsPOP.Range("C4").Value = 0.2157
    $wsPOP.Range("D4").Value = 0.1818
    $wsPOP.Range("B5").Value = 0.149
    $wsPOP.Range("C5").Value = 0.223
    $wsPOP.Range("D5").Value = 0.504
    $wsPOP.Range("B6").Value = -0.3139
    $wsPOP.Range("C6").Value = 0.8305
    $wsPOP.Range("D6").Value = 0.7054
    $wsPOP.Range("B7").Value = 0.4095
    $wsPOP.Range("C7").Value = 0.2231
    $wsPOP.Range("D7").Value = 0.0664
    $wsPOP.Range("B8").Value = 0.3838
    $wsPOP.Range("C8").Value = 0.2084
    $wsPOP.Range("D8").Value = 0.0656
    $wsPOP.Range("B9").Value = -0.2931
    $wsPOP.Range("C9").Value = 0.8827
    $wsPOP.Range("D9").Value = 0.7399
    $wsPOP.Range("B10").Value = 0.3606
    $wsPOP.Range("C10").Value = 0.3086
    $wsPOP.Range("D10").Value = 0.2426
    $wsPOP.Range("B11").Value = 0.2009
    $wsPOP.Range("C11").Value = 0.3074
    $wsPOP.Range("D11").Value = 0.5134
    $wsPOP.Range("B12").Value = 0.3958
    $wsPOP.Range("C12").Value = 0.4672
    $wsPOP.Range("D12").Value = 0.3969
    $wsPOP.Range("B13").Value = 0.5143
    $wsPOP.Range("C13").Value = 0.5235
    $wsPOP.Range("D13").Value = 0.3259
    $wsPOP.Range("B14").Value = 0.5504
    $wsPOP.Range("C14").Value = 0.5427
    $wsPOP.Range("D14").Value = 0.3105
    $wsPOP.Range("B15").Value = 0.5372
    $wsPOP.Range("C15").Value = 0.5597
    $wsPOP.Range("D15").Value = 0.3372
    $wsPOP.Range("B16").Value = 0.5724
    $wsPOP.Range("C16").Value = 0.5642
    $wsPOP.Range("D16").Value = 0.3103
    $wsPOP.Range("B17").Value = 0.6515
    $wsPOP.Range("C17").Value = 0.5603
    $wsPOP.Range("D17").Value = 0.2449
    $wsPOP.Range("B18").Value = 0.597
    $wsPOP.Range("C18").Value = 0.5576
    $wsPOP.Range("D18").Value = 0.2844
    $wsPOP.Range("B19").Value = 0.501
    $wsPOP.Range("C19").Value = 0.5561
    $wsPOP.Range("D19").Value = 0.3676
    $wsPOP.Range("B20").Value = 0.5085
    $wsPOP.Range("C20").Value = 0.5452
    $wsPOP.Range("D20").Value = 0.351
    $wsPOP.Range("B21").Value = 0.5858
    $wsPOP.Range("C21").Value = 0.5525
    $wsPOP.Range("D21").Value = 0.289
    $wsPOP.Range("B22").Value = 0.6748
    $wsPOP.Range("C22").Value = 0.6711
    $wsPOP.Range("D22").Value = 0.3146
    $wsPOP.Range("B23").Value = -1.3452
    $wsPOP.Range("C23").Value = 0.088

# --- Pesticide sheet updates ---
    $wsPest.Range("B2").Value = 10.8778
    $wsPest.Range("C2").Value = 12.7466
    $wsPest.Range("D2").Value = 0.3934
    $wsPest.Range("B3").Value = 0.4866
    $wsPest.Range("C3").Value = 0.3391
    $wsPest.Range("D3").Value = 0.1512
    $wsPest.Range("B6").Value = 0.3978
    $wsPest.Range("C6").Value = 0.3188
    $wsPest.Range("D6").Value = 0.212
    $wsPest.Range("B10").Value = -6.361
    $wsPest.Range("C10").Value = 10.8758
    $wsPest.Range("D10").Value = 0.5586
    $wsPest.Range("B11").Value = -2.9097
    $wsPest.Range("C11").Value = 2.9619
    $wsPest.Range("D11").Value = 0.3259
    $wsPest.Range("B12").Value = -5.582
    $wsPest.Range("C12").Value = 5.05
    $wsPest.Range("D12").Value = 0.269
    $wsPest.Range("B13").Value = -7.6948
    $wsPest.Range("C13").Value = 7.2608
    $wsPest.Range("D13").Value = 0.2892
    $wsPest.Range("B14").Value = -9.2431
    $wsPest.Range("C14").Value = 9.4241
    $wsPest.Range("D14").Value = 0.3267
    $wsPest.Range("B15").Value = -10.2268
    $wsPest.Range("C15").Value = 11.2768
    $wsPest.Range("D15").Value = 0.3645
    $wsPest.Range("B16").Value = -10.6459
    $wsPest.Range("C16").Value = 12.5
    $wsPest.Range("D16").Value = 0.3944
    $wsPest.Range("B17").Value = -10.5861
    $wsPest.Range("C17").Value = 12.8267
    $wsPest.Range("D17").Value = 0.4092
    $wsPest.Range("B18").Value = -9.9481
    $wsPest.Range("C18").Value = 12.7392
    $wsPest.Range("D18").Value = 0.4349
    $wsPest.Range("B19").Value = -9.8913
    $wsPest.Range("C19").Value = 12.7469
    $wsPest.Range("D19").Value = 0.4378
    $wsPest.Range("B20").Value = -9.545
    $wsPest.Range("C20").Value = 12.755
    $wsPest.Range("D20").Value = 0.4543
    $wsPest.Range("B21").Value = -9.8953
    $wsPest.Range("C21").Value = 12.8016
    $wsPest.Range("D21").Value = 0.4395
    $wsPest.Range("B22").Value = -10.2585
    $wsPest.Range("C22").Value = 13.2301
    $wsPest.Range("D22").Value = 0.4381
    $wsPest.Range("B23").Value = -0.6258
    $wsPest.Range("C23").Value = 0.1822
    $wsPest.Range("D23").Value = 0.0006

